$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Insert a new row at 13 to host the "Docentes responsaveis" data row
#    (moves old rows 13..23 down to 14..24)
$ws.Rows("13:13").Insert()

# Copy the body-row formatting (style 2 / style 3) from row 2 into the new row 13
$ws.Range("B2:C2").Copy()
$ws.Range("B13:C13").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Clear the inherited label cell in column A of the new row (no label on this row)
$ws.Range("A13").Clear()

# Fill the newly inserted row with the "responsible lecturer" data,
# relocated here from its previous (incorrect) position
$ws.Range("B13").Value = '5840751 - Jayne Carlos de Souza Barboza'
$ws.Range("C13").Value = '5840751 - Jayne Carlos de Souza Barboza'

# 2) Objetivos (row 10): replace the misplaced text with the real Portuguese objectives
$ws.Range("B10").Value = 'Introduzir o aluno ao conhecimento e entendimento das funções da química orgânica. Compreender as condições para que as reações orgânicas ocorram. Saber analisar os produtos das reações empregando espectroscopia na região do infravermelho e por ressonância magnética nuclear de próton e carbono. Levá-los a compreender os procedimentos e problemas industriais, sociais e ambientais, com os quais a Engenharia Química está estritamente relacionada, tornando-os, dessa forma, capazes a exercerem a função do Engenheiro Químico e realizar as mudanças que se fizerem necessárias.'
$ws.Range("C10").Value = 'Introduzir o aluno ao conhecimento e entendimento das funções da química orgânica. Compreender as condições para que as reações orgânicas ocorram. Saber analisar os produtos das reações empregando espectroscopia na região do infravermelho e por ressonância magnética nuclear de próton e carbono. Levá-los a compreender os procedimentos e problemas industriais, sociais e ambientais, com os quais a Engenharia Química está estritamente relacionada, tornando-os, dessa forma, capazes a exercerem a função do Engenheiro Químico e realizar as mudanças que se fizerem necessárias.'

# 3) Programa resumido (now row 14): replace placeholder with the Portuguese short syllabus
$ws.Range("B14").Value = 'Ácidos carboxílicos e seus derivados (haletos de acila, anidridos, ésteres e amidas – 16 horas/aulas). Aldeídos e cetonas (4 horas/aulas). Aminas (4 horas/aulas). Fenóis (3 horas/aulas). Sais de diazônio e azo compostos (3 horas/aulas). Ácidos sulfônicos e derivados (4 horas/aulas). Cetoácidos (3 horas/aulas). Hidroxi ácidos (3 aulas). Heterocíclicos (4 horas/aulas). Noções de síntese orgânica (8 horas/aulas). Análises por espectroscopia na região do infravermelho e por ressonância magnética nuclear de próton e carbono (8 horas/aulas)'
$ws.Range("C14").Value = 'Ácidos carboxílicos e seus derivados (haletos de acila, anidridos, ésteres e amidas – 16 horas/aulas). Aldeídos e cetonas (4 horas/aulas). Aminas (4 horas/aulas). Fenóis (3 horas/aulas). Sais de diazônio e azo compostos (3 horas/aulas). Ácidos sulfônicos e derivados (4 horas/aulas). Cetoácidos (3 horas/aulas). Hidroxi ácidos (3 aulas). Heterocíclicos (4 horas/aulas). Noções de síntese orgânica (8 horas/aulas). Análises por espectroscopia na região do infravermelho e por ressonância magnética nuclear de próton e carbono (8 horas/aulas)'

# 4) Programa (now row 16): replace placeholder with the Portuguese full syllabus
$ws.Range("B16").Value = 'Ácidos carboxílicos e derivados: Nomenclatura, propriedades físicas, processos de obtenção e propriedades químicas. Reações de substituição nucleofilícas em carbono acíclico. Síntese dos haletos de acila, anidridos, ésteres, amidas, aldeídos, cetonas, aminas, fenóis e suas propriedades químicas. Alfa halo-ácidos. Cetoácidos: processos de obtenção e propriedades químicas. Preparação de alfa cetoésteres. Hidroxi ácidos: processos de obtenção e propriedades químicas. Ácidos sulfônicos e seus derivados: processos de obtenção e propriedades químicas. Cloreto de sulfonila, sulfonamidas. Ésteres dos ácidos sulfônicos. Aminas: propriedades físicas, processos de obtenção e propriedades químicas. Reações de acoplamento dos sais de diazônio e azo compostos: Fenóis. Noções de síntese orgânica. Preparação de amostras, análises nos equipamentos e análises de espectros de compostos orgânicos na região do infravermelho e por ressonância magnética nuclear de próton e carbono'
$ws.Range("C16").Value = 'Ácidos carboxílicos e derivados: Nomenclatura, propriedades físicas, processos de obtenção e propriedades químicas. Reações de substituição nucleofilícas em carbono acíclico. Síntese dos haletos de acila, anidridos, ésteres, amidas, aldeídos, cetonas, aminas, fenóis e suas propriedades químicas. Alfa halo-ácidos. Cetoácidos: processos de obtenção e propriedades químicas. Preparação de alfa cetoésteres. Hidroxi ácidos: processos de obtenção e propriedades químicas. Ácidos sulfônicos e seus derivados: processos de obtenção e propriedades químicas. Cloreto de sulfonila, sulfonamidas. Ésteres dos ácidos sulfônicos. Aminas: propriedades físicas, processos de obtenção e propriedades químicas. Reações de acoplamento dos sais de diazônio e azo compostos: Fenóis. Noções de síntese orgânica. Preparação de amostras, análises nos equipamentos e análises de espectros de compostos orgânicos na região do infravermelho e por ressonância magnética nuclear de próton e carbono'

# 5) Metodo (now row 19): set the evaluation method text
$ws.Range("B19").Value = 'Duas provas semestrais teóricas (P1 e P2).'
$ws.Range("C19").Value = 'Duas provas semestrais teóricas (P1 e P2).'

# 6) Criterio (now row 20): set the grading criterion text
$ws.Range("B20").Value = 'A média final (M) será calculada pela expressão M = (P1 + 2 x P2)/3'
$ws.Range("C20").Value = 'A média final (M) será calculada pela expressão M = (P1 + 2 x P2)/3'

# 7) Norma de recuperacao (now row 21): set the recovery norm text
$ws.Range("B21").Value = 'Aos alunos que tiverem freqüência mínima de 70% e média final menor que 5,0 e igual ou maior que 3,0, será dada  recuperação  com uma avaliação escrita. A média dessa avaliação somada com a média anterior das P1 e P2, se superior a cinco (5,0), levará a aprovação do aluno.'
$ws.Range("C21").Value = 'Aos alunos que tiverem freqüência mínima de 70% e média final menor que 5,0 e igual ou maior que 3,0, será dada  recuperação  com uma avaliação escrita. A média dessa avaliação somada com a média anterior das P1 e P2, se superior a cinco (5,0), levará a aprovação do aluno.'

# 8) Bibliografia (now row 22): set the bibliography text
$ws.Range("B22").Value = '1) BARBOSA, L. C. A. Química Orgânica. Viçosa: Editora UFV, 2000. 2) BRESLOW, R. Questões e Exercícios de Química Orgânica. São Paulo: Makrons Books Editora, 1996. 3) CAMPOS, M. M.. Química Orgânica. São Paulo: Editora Prentice Hall, 2006. 4) DURST, H. D. Fundamentos de Química Orgânica. São Paulo: Editora Edgard Blucher, 1997. 5) HENDRIKSON, J. B.; CRAM, D. J. Mecanismos de Reações Orgânicas. São Paulo: Livraria Editora, 1966. 6) MCMURRY, J. Química Orgânica. São Paulo: Editora Pioneira Thomson Leraning, 2005. 7) SOLOMONS, T.W.G; FRYHLE, G.. Química Orgânica. Rio de Janeiro: Livros Técnicos e Científicos Editora, 2001. 8) SOARES, B. G. et al. Química Orgânica Experimental. Barcelona: Editorial Reverte, 1985. 9) Silverstein, R. M.; Bassler, G. C.; Morrill, T. C. Identificação Espectrométrica de Compostos Orgânicos. Guanabara Koogan'
$ws.Range("C22").Value = '1) BARBOSA, L. C. A. Química Orgânica. Viçosa: Editora UFV, 2000. 2) BRESLOW, R. Questões e Exercícios de Química Orgânica. São Paulo: Makrons Books Editora, 1996. 3) CAMPOS, M. M.. Química Orgânica. São Paulo: Editora Prentice Hall, 2006. 4) DURST, H. D. Fundamentos de Química Orgânica. São Paulo: Editora Edgard Blucher, 1997. 5) HENDRIKSON, J. B.; CRAM, D. J. Mecanismos de Reações Orgânicas. São Paulo: Livraria Editora, 1966. 6) MCMURRY, J. Química Orgânica. São Paulo: Editora Pioneira Thomson Leraning, 2005. 7) SOLOMONS, T.W.G; FRYHLE, G.. Química Orgânica. Rio de Janeiro: Livros Técnicos e Científicos Editora, 2001. 8) SOARES, B. G. et al. Química Orgânica Experimental. Barcelona: Editorial Reverte, 1985. 9) Silverstein, R. M.; Bassler, G. C.; Morrill, T. C. Identificação Espectrométrica de Compostos Orgânicos. Guanabara Koogan'
